$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "2005.2. MONTGOMERY",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2005.^l2. MONTGOMERY", 2)

$d.Content.Find.Execute(
    "2004.3. GRANT",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2004.^l3. GRANT", 2)

$d.Content.Find.Execute(
    "1996.4. WERKENA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1996.^l4. WERKENA", 2)
